# Commit: "swap" the presentation's applied theme ("Integral") for the
# stock "Office Theme" colors (and vice-versa for the theme used by the
# notes master) - as captured by the ppt/theme/theme1.xml <-> theme2.xml
# content exchange in the OOXML diff.
#
# The slide master (and therefore the whole deck's visible design) is
# backed by ppt/theme/theme2.xml; we recolor it, through the standard
# ColorScheme object, to the 12 "Office Theme" RGB values so that the
# master's theme becomes the Office palette.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

# PpColorSchemeIndex order matches a:clrScheme child order:
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$cs.Colors(1).RGB  = 0          # dk1      #000000
$cs.Colors(2).RGB  = 16777215   # lt1      #FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      #44546A
$cs.Colors(4).RGB  = 15132391   # lt2      #E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  #5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  #ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  #A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  #FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  #4472C4
$cs.Colors(10).RGB = 4697456    # accent6  #70AD47
$cs.Colors(11).RGB = 12673797   # hlink    #0563C1
$cs.Colors(12).RGB = 7491477    # folHlink #954F72
